$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Events 2026"

# ---------------------------------------------------------------------------
# New event rows (7-14) for the "Events 2026" sheet.
# Columns: A=event_title, B=date, C=time, D=venue, E=type,
#          F=description, G=collaborators, H=catering, I=signup_link
# ---------------------------------------------------------------------------

# --- Row 7 : Poker Workshop -------------------------------------------------
$ws.Cells.Item(7,1).Value = "Poker Workshop"
$ws.Cells.Item(7,2).Value = "2026-03-04"
$ws.Cells.Item(7,2).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(7,3).Value = "TBA"
$ws.Cells.Item(7,4).Value = "TBA"
$ws.Cells.Item(7,5).Value = "social x academic"
$ws.Cells.Item(7,6).Value = "Crash course in poker ahead of our poker tournament "
$ws.Cells.Item(7,6).Style = "Normal"
$ws.Cells.Item(7,7).Value = "SUDATA, SQT"
$ws.Cells.Item(7,8).Value = "TBA"
$ws.Cells.Item(7,8).Style = "Normal"
$ws.Cells.Item(7,9).Value = "TBA"

# Propagate the date style (created above) to the other new-row date cells
# so that every B-cell below shares the very same style entry.
$ws.Cells.Item(7,2).Copy($ws.Range("B8:B12"))

# --- Row 8 : Introduction to R Workshop ------------------------------------
$ws.Cells.Item(8,1).Value = "Introduction to R Workshop"
$ws.Cells.Item(8,2).Value = "2026-03-06"
$ws.Cells.Item(8,3).Value = 0.66666666666666663
$ws.Cells.Item(8,3).NumberFormat = "h:mm"
$ws.Cells.Item(8,4).Value = "TBA"
$ws.Cells.Item(8,5).Value = "academic"
$ws.Cells.Item(8,6).Value = "Learn the fundamentals of R programming for data analysis, and connect with new friends in a relaxed, supportive environment"
$ws.Cells.Item(8,6).Style = "Normal"
$ws.Cells.Item(8,7).Value = "SUDATA"
$ws.Cells.Item(8,8).Value = "Snacks"
$ws.Cells.Item(8,8).Style = "Normal"
$ws.Cells.Item(8,9).Value = "TBA"

# --- Row 9 : Meet The Startups ----------------------------------------------
$ws.Cells.Item(9,1).Value = "Meet The Startups"
$ws.Cells.Item(9,2).Value = "2026-03-16"
$ws.Cells.Item(9,3).Value = "TBA"
$ws.Cells.Item(9,4).Value = "TBA"
$ws.Cells.Item(9,5).Value = "academic"
$ws.Cells.Item(9,6).Value = "Panel discussion with startup founders about data science in entrepreneurship."
$ws.Cells.Item(9,6).Style = "Normal"
$ws.Cells.Item(9,7).Value = "SUDATA, SUSMI, STARTUPLINK, MEDSCISOC"
$ws.Cells.Item(9,8).Value = "TBA"
$ws.Cells.Item(9,8).Style = "Normal"
$ws.Cells.Item(9,9).Value = "TBA"

# --- Row 10 : Datathon -------------------------------------------------------
$ws.Cells.Item(10,1).Value = "Datathon"
$ws.Cells.Item(10,2).Value = "2026-03-30"
$ws.Cells.Item(10,3).Value = "TBA"
$ws.Cells.Item(10,4).Value = "TBA"
$ws.Cells.Item(10,5).Value = "academic"
$ws.Cells.Item(10,6).Value = "Multi-day datathon competition. Solve real-world data problems and win prizes!"
$ws.Cells.Item(10,6).Style = "Normal"
$ws.Cells.Item(10,7).Value = "SUDATA, COMM-STEM, SYNCS"
$ws.Cells.Item(10,8).Value = "Lunch and Snacks"
$ws.Cells.Item(10,8).Style = "Normal"
$ws.Cells.Item(10,9).Value = "TBA"

# --- Row 11 : Amstelveen Consulting Workshop --------------------------------
$ws.Cells.Item(11,1).Value = "Amstelveen Consulting Workshop "
$ws.Cells.Item(11,2).Value = "2026-04-20"
$ws.Cells.Item(11,3).Value = "TBA"
$ws.Cells.Item(11,4).Value = "TBA"
$ws.Cells.Item(11,5).Value = "academic"
$ws.Cells.Item(11,6).Value = "Consulting workshop and short case competition practice in collaboration with Amstelveen"
$ws.Cells.Item(11,6).Style = "Normal"
$ws.Cells.Item(11,7).Value = "SUDATA, BISA"
$ws.Cells.Item(11,8).Value = "TBA"
$ws.Cells.Item(11,8).Style = "Normal"
$ws.Cells.Item(11,9).Value = "TBA"

# --- Row 12 : Data Science Careers Panel ------------------------------------
$ws.Cells.Item(12,1).Value = "Data Science Careers Panel"
$ws.Cells.Item(12,2).Value = "2026-05-04"
$ws.Cells.Item(12,3).Value = "TBA"
$ws.Cells.Item(12,4).Value = "TBA"
$ws.Cells.Item(12,5).Value = "academic"
$ws.Cells.Item(12,6).Value = "Learn about data science fundamentals and hear from industry professionals about career paths."
$ws.Cells.Item(12,6).Style = "Normal"
$ws.Cells.Item(12,7).Value = "SUDATA, 180 Degrees"
$ws.Cells.Item(12,8).Value = "TBA"
$ws.Cells.Item(12,8).Style = "Normal"
$ws.Cells.Item(12,9).Value = "TBA"

# --- Row 13 : Industry Networking Night (explicit black font) --------------
$ws.Cells.Item(13,1).Value = "Industry Networking Night"
$ws.Cells.Item(13,1).Font.Color = 0
# Propagate the plain black-font style to the rest of row 13 and row 14
# (column B handled separately below since it also carries a date format).
$ws.Cells.Item(13,1).Copy($ws.Range("C13:I13"))
$ws.Cells.Item(13,1).Copy($ws.Range("A14:A14"))
$ws.Cells.Item(13,1).Copy($ws.Range("C14:I14"))

$ws.Cells.Item(13,2).Value = "2026-05-18"
$ws.Cells.Item(13,2).Font.Color = 0
$ws.Cells.Item(13,2).NumberFormat = "mm-dd-yy"
$ws.Cells.Item(13,2).Copy($ws.Range("B14:B14"))

$ws.Cells.Item(13,3).Value = "TBA"
$ws.Cells.Item(13,4).Value = "TBA"
$ws.Cells.Item(13,5).Value = "academic"
$ws.Cells.Item(13,6).Value = "Major networking event with industry professionals from various data science companies."
$ws.Cells.Item(13,7).Value = "SUDATA"
$ws.Cells.Item(13,8).Value = "Canapés and drinks"
$ws.Cells.Item(13,9).Value = "TBA"

# --- Row 14 : DATA1001/2001 Revision Session (explicit black font) ---------
$ws.Cells.Item(14,1).Value = "DATA1001/2001 Revision Session"
$ws.Cells.Item(14,2).Value = "2026-06-01"
$ws.Cells.Item(14,3).Value = "TBA"
$ws.Cells.Item(14,4).Value = "TBA"
$ws.Cells.Item(14,5).Value = "academic"
$ws.Cells.Item(14,6).Value = "STUVAC revision session for DATA1001 and DATA2001 final exam preparation."
$ws.Cells.Item(14,7).Value = "SUDATA"
$ws.Cells.Item(14,8).Value = "Lunch and Snacks"
$ws.Cells.Item(14,9).Value = "TBA"

# Update the visible selection to mirror the authored workbook.
$ws.Range("D23").Select()
